# fix bug display error when import excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (student 1): roll number / name were shifted/garbled on import - fix them,
# and clear the stray Nationality value that didn't belong on this row.
$ws.Range("B2").Value = "HE130576"
$ws.Range("C2").Value = "Phạm Thanh Hà0"
$ws.Range("G2").ClearContents()

# Row 3 (student 2): RollNumber/Fullname/PlaceOfBirth were duplicated from
# other rows by mistake - clear them, and set the correct Content value.
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("K3").Value = "Tình một đêm"

# Row 4 (student 3): same stray duplicate fix, plus correcting Content text.
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("K4").Value = "Sáo"

# Row 5 (student 4): clear stray duplicate Fullname/Nationality, correct Content text.
$ws.Range("C5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("K5").Value = "Lừa trái tim đàn bà"

# Restore the active selection to C2, matching the saved workbook state.
$ws.Range("C2").Select()
